# Apply "想去人数" (want-to-go count) updates to the "展览" and "全部类型"
# sheets, matching the data refresh captured by the diff.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 85
$ws1.Range("F6").Value = 10078
$ws1.Range("F8").Value = 912
$ws1.Range("F9").Value = 1253
$ws1.Range("F10").Value = 6063
$ws1.Range("F12").Value = 359
$ws1.Range("F13").Value = 185
$ws1.Range("F15").Value = 3092
$ws1.Range("F17").Value = 300
$ws1.Range("F18").Value = 596
$ws1.Range("F23").Value = 1537

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 85
$ws4.Range("F7").Value = 10078
$ws4.Range("F9").Value = 912
$ws4.Range("F10").Value = 1253
$ws4.Range("F11").Value = 6063
$ws4.Range("F13").Value = 359
$ws4.Range("F14").Value = 185
$ws4.Range("F16").Value = 3092
$ws4.Range("F18").Value = 300
$ws4.Range("F19").Value = 596
$ws4.Range("F24").Value = 1537
